$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row11 = @(10, "Staff", 1, 0, 8.5, 7.8, 14, 6.3, 9.2, 5.7, 9.4, 3.2, 7.8, 4.4, 4.8, 3.6, 7.5, 8.3, 5.1, 5.7, 0)
$row12 = @(11, "Staff", 1, 0, 15.5, 2.9, 14.3, 12, 10.2, 1.1, 14.4, 1.6, 1.7, 3.3, 3.7, 3.8, 12.8, 2.3, 2, 5.6, 0)

for ($c = 0; $c -lt $row11.Length; $c++) {
    $ws.Cells.Item(11, $c + 1).Value = $row11[$c]
}

for ($c = 0; $c -lt $row12.Length; $c++) {
    $ws.Cells.Item(12, $c + 1).Value = $row12[$c]
}

$ws.Range("U12").Select()
